$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.457.42"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.853.41"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'233.29"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4744"
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("D8").Value = "'0.2762"
$ws.Range("E8").Value = "  +3.10%  "
$ws.Range("D9").Value = "'0.06346"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'18.02"
$ws.Range("E10").Value = "  +12.18%  "
$ws.Range("D11").Value = "1.903.88"
$ws.Range("E11").Value = "  +4.09%  "
$ws.Range("D12").Value = "'0.07465"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "'4.988"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "'84.78"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "'0.6246"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "30.424.85"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "'246.60"
$ws.Range("E17").Value = "  +8.93%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "'12.71"
$ws.Range("E19").Value = "  +2.68%  "
$ws.Range("D20").Value = "'0.000007350"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").Value = "'0.9994"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'4.934"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").Value = "'5.907"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "'164.06"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").Value = "'9.012"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "'18.00"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("D27").Value = "'1.880"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").Value = "'0.1025"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").Value = "'4.043"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").Value = "'3.839"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("D32").Value = "'0.04832"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "'1.131"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'0.6983"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").Value = "'2.694"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'0.01891"
$ws.Range("E36").Value = "  +4.11%  "
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("D38").Value = "'0.8760"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "'1.984"
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("D40").Value = "'106.32"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'0.4067"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").Value = "'5.506"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").Value = "'7.187"
$ws.Range("E44").Value = "  +3.14%  "
$ws.Range("D45").Value = "'63.40"
$ws.Range("E45").Value = "  +6.31%  "
$ws.Range("D46").Value = "'0.1200"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "'34.07"
$ws.Range("E47").Value = "  +4.11%  "
$ws.Range("D48").Value = "'8.536"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").Value = "'0.05496"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "'1.352"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "'0.3689"
$ws.Range("E51").Value = "  +1.85%  "
